$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns F and G in the header row, and new data (0) values in rows 2-4
$ws.Range("F1").Value = "BL "
$ws.Range("G1").Value = "Operating frequency"

$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Header row (A1:E1) and C1 lose their previously-applied center alignment / format
$ws.Range("A1:E1").ClearFormats()

# Move the active selection to E10
$null = $ws.Range("E10").Select()
